$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 574; this shifts existing rows 574:676 down to 575:677
$ws.Rows(574).Insert()

# Populate the newly inserted row 574 with the new record
$ws.Range("A574").Value = 4
$ws.Range("B574").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C574").Value = 'Los Lagos'
$ws.Range("D574").Value = 45180
$ws.Range("E574").Value = 10
$ws.Range("F574").Value = 'Fruta'
$ws.Range("G574").Value = 100102
$ws.Range("H574").Value = 'Cítricos'
$ws.Range("I574").Value = 100102006
$ws.Range("J574").Value = 'Pomelo'
$ws.Range("K574").Value = 'Start Ruby'
$ws.Range("L574").Value = 'Primera'
$ws.Range("M574").Value = 30
$ws.Range("N574").Value = 15000
$ws.Range("O574").Value = 16000
$ws.Range("P574").Value = 15500
$ws.Range("Q574").Value = '$/caja 14 kilos empedrada'
$ws.Range("R574").Value = "Región de O'Higgins"
$ws.Range("S574").Value = 1107
$ws.Range("T574").Value = 14
